$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '53.117.81'
$ws.Range("E2").Value = '  -10.23%  '
$ws.Range("D3").Value = '2.379.14'
$ws.Range("E3").Value = '  -13.09%  '
Set-TextValue $ws "D4" '0.999'
$ws.Range("E4").Value = '  -0.02%  '
Set-TextValue $ws "D5" '457.43'
$ws.Range("E5").Value = '  -10.20%  '
Set-TextValue $ws "D6" '128.93'
$ws.Range("E6").Value = '  -9.15%  '
Set-TextValue $ws "D7" '0.995'
$ws.Range("E7").Value = '  -0.25%  '
Set-TextValue $ws "D8" '0.479'
$ws.Range("E8").Value = '  -10.51%  '
$ws.Range("D9").Value = '2.399.06'
$ws.Range("E9").Value = '  -12.74%  '
Set-TextValue $ws "D10" '0.0939'
$ws.Range("E10").Value = '  -10.55%  '
Set-TextValue $ws "D11" '5.23'
$ws.Range("E11").Value = '  -14.95%  '
Set-TextValue $ws "D12" '0.310'
$ws.Range("E12").Value = '  -11.57%  '
$ws.Range("E13").Value = '  -4.70%  '
$ws.Range("D14").Value = '2.788.60'
$ws.Range("E14").Value = '  -13.21%  '
$ws.Range("D15").Value = '53.135.07'
$ws.Range("E15").Value = '  -9.94%  '
Set-TextValue $ws "D16" '19.41'
$ws.Range("E16").Value = '  -11.41%  '
Set-TextValue $ws "D17" '0.0000129'
$ws.Range("E17").Value = '  -5.80%  '
$ws.Range("D18").Value = '2.370.95'
$ws.Range("E18").Value = '  -13.02%  '
Set-TextValue $ws "D19" '4.12'
$ws.Range("E19").Value = '  -13.57%  '
Set-TextValue $ws "D20" '304.97'
$ws.Range("E20").Value = '  -12.40%  '
Set-TextValue $ws "D21" '9.28'
$ws.Range("E21").Value = '  -16.10%  '
Set-TextValue $ws "D22" '0.999'
$ws.Range("E22").Value = '  -0.35%  '
$ws.Range("E23").Value = '  +0.50%  '
Set-TextValue $ws "D24" '5.28'
$ws.Range("E24").Value = '  -16.05%  '
Set-TextValue $ws "D25" '55.30'
$ws.Range("E25").Value = '  -12.36%  '
Set-TextValue $ws "D26" '1.00'
$ws.Range("E26").Value = '  +1.03%  '
Set-TextValue $ws "D27" '0.380'
$ws.Range("E27").Value = '  -10.97%  '
$ws.Range("D28").Value = '2.466.20'
$ws.Range("E28").Value = '  -13.07%  '
Set-TextValue $ws "D29" '0.150'
$ws.Range("E29").Value = '  -13.36%  '
Set-TextValue $ws "D30" '7.00'
$ws.Range("E30").Value = '  -7.00%  '
Set-TextValue $ws "D31" '0.995'
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").Value = '0.0₃0718'
$ws.Range("E32").Value = '  -14.67%  '
Set-TextValue $ws "D33" '144.40'
$ws.Range("E33").Value = '  -3.32%  '
Set-TextValue $ws "D34" '17.50'
$ws.Range("E34").Value = '  -8.91%  '
Set-TextValue $ws "D35" '1.41'
$ws.Range("E35").Value = '  -13.33%  '
Set-TextValue $ws "D36" '4.94'
$ws.Range("E36").Value = '  -8.77%  '
Set-TextValue $ws "D37" '3.48'
$ws.Range("E37").Value = '  -18.11%  '
Set-TextValue $ws "D38" '1.05'
$ws.Range("E38").Value = '  -8.74%  '
Set-TextValue $ws "D39" '0.791'
$ws.Range("E39").Value = '  -17.05%  '
Set-TextValue $ws "D40" '0.995'
$ws.Range("E40").Value = '  +0.02%  '
Set-TextValue $ws "D41" '32.91'
$ws.Range("E41").Value = '  -9.19%  '
Set-TextValue $ws "D42" '0.588'
$ws.Range("E42").Value = '  -2.89%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws "D43" '3.23'
$ws.Range("E43").Value = '  -9.05%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws "D44" '0.0519'
$ws.Range("E44").Value = '  -7.75%  '
Set-TextValue $ws "D45" '10.09'
$ws.Range("E45").Value = '  -2.56%  '
Set-TextValue $ws "D46" '1.22'
$ws.Range("E46").Value = '  -12.40%  '
$ws.Range("D47").Value = '1.920.53'
$ws.Range("E47").Value = '  -12.42%  '
Set-TextValue $ws "D48" '0.0215'
$ws.Range("E48").Value = '  -6.72%  '
Set-TextValue $ws "D49" '0.0858'
$ws.Range("E49").Value = '  -3.47%  '
Set-TextValue $ws "D50" '4.15'
$ws.Range("E50").Value = '  -14.03%  '
Set-TextValue $ws "D51" '16.32'
$ws.Range("E51").Value = '  -15.03%  '
